$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data values as per latest scrape.

$ws.Range("D2").Value = '68.013.23'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '3.270.19'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.75'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.96'
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.604'
$ws.Range("E8").Value = '  +1.26%  '
$ws.Range("E9").Value = '  -2.27%  '
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  -3.55%  '
$ws.Range("D12").Value = '3.834.92'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.139'
$ws.Range("E13").Value = '  +0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.43'
$ws.Range("E14").Value = '  -3.82%  '
$ws.Range("D15").Value = '68.029.30'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("D17").Value = '3.286.36'
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.42'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '408.74'
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.55'
$ws.Range("E21").Value = '  -1.68%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.15'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.510'
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000118'
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.188'
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.47'
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.94'
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.71'
$ws.Range("E30").Value = '  -1.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.46'
$ws.Range("E31").Value = '  -4.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.91'
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -2.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.29'
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("E36").Value = '  -3.42%  '
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.27'
$ws.Range("E38").Value = '  +3.45%  '
$ws.Range("E39").Value = '  -3.36%  '
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("E41").Value = '  -3.55%  '
$ws.Range("D42").Value = '2.663.27'
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0677'
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.66'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '334.63'
$ws.Range("E47").Value = '  -3.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0274'
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.31'
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.971'
$ws.Range("E51").Value = '  -1.17%  '
